$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors used by the two new small "helper list" fonts (Arial 8pt).
$gray = 4473924   # RGB(0x44,0x44,0x44) -> FF444444
$black = 0        # RGB(0,0,0)          -> FF000000

function SetCell($addr, $val, $color) {
  $c = $ws.Range($addr)
  $c.Value = $val
  $c.Font.Name = "Arial"
  $c.Font.Size = 8
  $c.Font.Color = $color
}

# ---- Column Z: "Jenis Perangkat" helper list (used by the B2:B dropdown) ----
$ws.Range("Z1").Value = "Jenis Perangkat"
SetCell "Z2"  "Laptop"    $gray
SetCell "Z3"  "Komputer"  $black
SetCell "Z4"  "Handphone" $black
SetCell "Z5"  "Hardware"  $black
SetCell "Z6"  "Software"  $black
SetCell "Z7"  "Printer"   $gray
SetCell "Z8"  "Kamera"    $black
SetCell "Z9"  "Mouse"     $black
SetCell "Z10" "Parabola"  $black
SetCell "Z11" "Kabel"     $black

# ---- Column AA: "Status Kebutuhan" helper list (used by the I2:I dropdown) ----
$ws.Range("AA1").Value = "Status Kebutuhan"
SetCell "AA2" "Pergantian Barang" $black
SetCell "AA3" "Perbaikan Barang"  $gray
SetCell "AA4" "Kerusakan Barang"  $gray
$ws.Range("AA5").Value = "Permintaan Barang"

# ---- Column AB: "Department" helper list (used by the K2:K dropdown) ----
$ws.Range("AB1").Value  = "Department"
$ws.Range("AB2").Value  = "Board Of Director"
$ws.Range("AB3").Value  = "Billing Support"
$ws.Range("AB4").Value  = "Banking"
$ws.Range("AB5").Value  = "Finance & Accounting"
$ws.Range("AB6").Value  = "Human Resource"
$ws.Range("AB7").Value  = "HUB Operation"
$ws.Range("AB8").Value  = "Legal"
$ws.Range("AB9").Value  = "MP Upgrade"
$ws.Range("AB10").Value = "General Affair"
$ws.Range("AB11").Value = "Services Delivery"
$ws.Range("AB12").Value = "Product Development"
$ws.Range("AB13").Value = "Purchasing"
$ws.Range("AB14").Value = "QMR"
$ws.Range("AB15").Value = "Sales & Marketing"
$ws.Range("AB16").Value = "Services"
$ws.Range("AB17").Value = "Warehouse & Logistic"
$ws.Range("AB18").Value = "Workshop"
$ws.Range("AB19").Value = "Business Support"
$ws.Range("AB20").Value = "NIX"
$ws.Range("AB21").Value = "Bitnet"

# ---- Data validation dropdowns driven by the helper lists above ----
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, "=`$Z`$2:`$Z`$11")
$ws.Range("I2:I1048576").Validation.Add(3, 1, 1, "=`$AA`$2:`$AA`$5")
$ws.Range("K2:K1048576").Validation.Add(3, 1, 1, "=`$AB`$2:`$AB`$21")

# ---- Misc cosmetic adjustments ----
$ws.Columns.Item(9).ColumnWidth = 15.49869791666667
$ws.Columns.Item(26).ColumnWidth = 13.385416666666666
$ws.Columns.Item(27).ColumnWidth = 15.944010416666666

$ws.PageSetup.PaperSize = 9

$ws.Range("E18").Select()
